$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 8 (215f1916 file) - newly stamped handoff datetimes
$zhcn.Range("D8").Value = "2016-03-09 10:30:32"
$dede.Range("D8").Value = "2016-03-09 10:30:43"

# Row 10 (b626bd2b file) - reuses the same new handoff datetimes
$zhcn.Range("D10").Value = "2016-03-09 10:30:32"
$dede.Range("D10").Value = "2016-03-09 10:30:43"
